$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-15 21:18:32'
$ws.Range("E3").Value = '2026-02-15 21:18:35'
$ws.Range("E4").Value = '2026-02-15 21:18:37'
$ws.Range("H4").Value = '72%'
$ws.Range("E5").Value = '2026-02-15 21:18:40'
$ws.Range("I5").Value = '7.3 mm'
$ws.Range("O5").Value = '-4.4 °C'
$ws.Range("E6").Value = '2026-02-15 21:18:42'
$ws.Range("H6").Value = '61%'
$ws.Range("E7").Value = '2026-02-15 21:18:45'
$ws.Range("O7").Value = '11.9 °C'
$ws.Range("E8").Value = '2026-02-15 21:18:47'
$ws.Range("O8").Value = '8.3 °C'
$ws.Range("E9").Value = '2026-02-15 21:18:50'
$ws.Range("H9").Value = '52%'
$ws.Range("O9").Value = '10.8 °C'
$ws.Range("E10").Value = '2026-02-15 21:18:52'
$ws.Range("K10").Value = '11.8 MJ/m2'
$ws.Range("E11").Value = '2026-02-15 21:18:55'
$ws.Range("H11").Value = '45%'
$ws.Range("O11").Value = '7.1 °C'
$ws.Range("E12").Value = '2026-02-15 21:18:57'
$ws.Range("H12").Value = '58%'
$ws.Range("E13").Value = '2026-02-15 21:19:00'
$ws.Range("H13").Value = '38%'
$ws.Range("E14").Value = '2026-02-15 21:19:02'
$ws.Range("H14").Value = '59%'
$ws.Range("O14").Value = '10.7 °C'
$ws.Range("E15").Value = '2026-02-15 21:19:04'
$ws.Range("H15").Value = '52%'
$ws.Range("O15").Value = '10.5 °C'
$ws.Range("E16").Value = '2026-02-15 21:19:07'
$ws.Range("I16").Value = '1.1 mm'
$ws.Range("O16").Value = '-1.9 °C'
$ws.Range("E17").Value = '2026-02-15 21:19:09'
$ws.Range("H17").Value = '40%'
$ws.Range("E18").Value = '2026-02-15 21:19:12'
$ws.Range("O18").Value = '7.6 °C'
$ws.Range("E19").Value = '2026-02-15 21:19:14'
$ws.Range("O19").Value = '3.5 °C'
$ws.Range("E20").Value = '2026-02-15 21:19:17'
$ws.Range("H20").Value = '62%'
$ws.Range("L20").Value = '79.2 km/h - 330º 20:33 TU'
$ws.Range("O20").Value = '-2.6 °C'
$ws.Range("E21").Value = '2026-02-15 21:19:19'
$ws.Range("E22").Value = '2026-02-15 21:19:22'
$ws.Range("E23").Value = '2026-02-15 21:19:24'
$ws.Range("I23").Value = '4.4 mm'
$ws.Range("O23").Value = '-3.5 °C'
$ws.Range("E24").Value = '2026-02-15 21:19:27'
$ws.Range("E25").Value = '2026-02-15 21:19:30'
$ws.Range("O25").Value = '-1.4 °C'
$ws.Range("E26").Value = '2026-02-15 21:19:32'
$ws.Range("E27").Value = '2026-02-15 21:19:34'
$ws.Range("H27").Value = '50%'
$ws.Range("E28").Value = '2026-02-15 21:19:37'
$ws.Range("H28").Value = '59%'
$ws.Range("J28").Value = '1015.7 hPa'
$ws.Range("E29").Value = '2026-02-15 21:19:39'
$ws.Range("H29").Value = '59%'
$ws.Range("K29").Value = '12.1 MJ/m2'
$ws.Range("O29").Value = '10.0 °C'
$ws.Range("E30").Value = '2026-02-15 21:19:42'
$ws.Range("H30").Value = '57%'
$ws.Range("E31").Value = '2026-02-15 21:19:44'
$ws.Range("O31").Value = '10.1 °C'
$ws.Range("E32").Value = '2026-02-15 21:19:47'
$ws.Range("O32").Value = '3.9 °C'
$ws.Range("E33").Value = '2026-02-15 21:19:49'
$ws.Range("E34").Value = '2026-02-15 21:19:52'
$ws.Range("L34").Value = '61.2 km/h - 49º 20:45 TU'
$ws.Range("E35").Value = '2026-02-15 21:19:54'
$ws.Range("E36").Value = '2026-02-15 21:19:57'
$ws.Range("E37").Value = '2026-02-15 21:20:00'
$ws.Range("H37").Value = '55%'
$ws.Range("E38").Value = '2026-02-15 21:20:02'
$ws.Range("O38").Value = '7.7 °C'
$ws.Range("E39").Value = '2026-02-15 21:20:05'
$ws.Range("E40").Value = '2026-02-15 21:20:07'
$ws.Range("H40").Value = '40%'
$ws.Range("J40").Value = '1016.2 hPa'
$ws.Range("O40").Value = '8.7 °C'
$ws.Range("E41").Value = '2026-02-15 21:20:09'
$ws.Range("O41").Value = '12.6 °C'
$ws.Range("E42").Value = '2026-02-15 21:20:12'
$ws.Range("H42").Value = '58%'
$ws.Range("E43").Value = '2026-02-15 21:20:14'
$ws.Range("E44").Value = '2026-02-15 21:20:17'
$ws.Range("H44").Value = '78%'
$ws.Range("I44").Value = '5.0 mm'
$ws.Range("E45").Value = '2026-02-15 21:20:19'
$ws.Range("I45").Value = '3.4 mm'
$ws.Range("E46").Value = '2026-02-15 21:20:22'
